$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for cells whose content could otherwise be
# auto-converted to a number/date by Excel (preserves literal formatting,
# e.g. "10.40" keeping its trailing zero, "0.998" etc).
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell $ws.Range("D2") "58.915.87"
$ws.Range("E2").Value = "  -2.35%  "
Set-TextCell $ws.Range("D3") "2.660.85"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextCell $ws.Range("D5") "523.94"
$ws.Range("E5").Value = "  -0.04%  "
Set-TextCell $ws.Range("D6") "144.31"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  +0.20%  "
Set-TextCell $ws.Range("D8") "0.569"
$ws.Range("E8").Value = "  -1.24%  "
Set-TextCell $ws.Range("D9") "6.96"
$ws.Range("E9").Value = "  +7.10%  "
$ws.Range("E10").Value = "  -3.49%  "
Set-TextCell $ws.Range("D11") "0.336"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  +1.41%  "
Set-TextCell $ws.Range("D13") "3.128.24"
$ws.Range("E13").Value = "  -1.27%  "
Set-TextCell $ws.Range("D14") "58.880.22"
$ws.Range("E14").Value = "  -2.47%  "
Set-TextCell $ws.Range("D15") "21.03"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("E16").Value = "  -1.93%  "
Set-TextCell $ws.Range("D17") "2.653.93"
$ws.Range("E17").Value = "  -6.70%  "
Set-TextCell $ws.Range("D18") "339.39"
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("E19").Value = "  -3.41%  "
Set-TextCell $ws.Range("D20") "10.40"
$ws.Range("E20").Value = "  -1.74%  "
Set-TextCell $ws.Range("D21") "6.36"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +0.26%  "
Set-TextCell $ws.Range("D23") "64.48"
$ws.Range("E23").Value = "  +2.27%  "
Set-TextCell $ws.Range("D24") "0.421"
$ws.Range("E24").Value = "  -0.43%  "
Set-TextCell $ws.Range("D25") "0.167"
$ws.Range("E25").Value = "  -1.23%  "
Set-TextCell $ws.Range("D26") "0.998"
$ws.Range("E26").Value = "  +0.41%  "
Set-TextCell $ws.Range("D27") "0.0₃0804"
$ws.Range("E27").Value = "  -2.03%  "
Set-TextCell $ws.Range("D28") "7.16"
$ws.Range("E28").Value = "  -2.89%  "
Set-TextCell $ws.Range("D29") "6.70"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("E30").Value = "  +0.10%  "
Set-TextCell $ws.Range("D31") "1.60"
$ws.Range("E31").Value = "  -0.26%  "
Set-TextCell $ws.Range("D32") "18.90"
$ws.Range("E32").Value = "  -1.40%  "
Set-TextCell $ws.Range("D33") "150.68"
$ws.Range("E34").Value = "  -3.35%  "
Set-TextCell $ws.Range("D35") "0.928"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("E36").Value = "  -4.97%  "
Set-TextCell $ws.Range("D37") "0.872"
$ws.Range("E37").Value = "  -0.79%  "
Set-TextCell $ws.Range("D38") "36.99"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D41") "0.614"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  +0.31%  "
Set-TextCell $ws.Range("D43") "275.64"
$ws.Range("E43").Value = "  -3.75%  "
Set-TextCell $ws.Range("D44") "19.73"
$ws.Range("E44").Value = "  -2.01%  "
Set-TextCell $ws.Range("D45") "0.0968"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("E46").Value = "  +1.97%  "
Set-TextCell $ws.Range("D47") "2.061.17"
$ws.Range("E47").Value = "  -3.63%  "
Set-TextCell $ws.Range("D48") "0.0535"
$ws.Range("E48").Value = "  -1.23%  "
Set-TextCell $ws.Range("D49") "4.74"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("E50").Value = "  -2.87%  "
Set-TextCell $ws.Range("D51") "18.84"
$ws.Range("E51").Value = "  -2.45%  "
